# Seeds.xlsx refactor: creating seeds from excel file
# - rename Products header "product_cathegory" -> "product_cathegory_id"
# - mark one product ("Tynk akrylowy", row 6) as promoted by filling
#   promoted_from / promoted_to with a timestamp
# - backfill a missing product_cathegory value for "Syropian fundamentowy 16 cm" (row 13)
# - add default "{}" avatars JSON for the seeded Users rows
# - drop the stray empty styled cell on ProductsCathegories!B1
# - widen the two new Products columns (F, G)

$wb = $excel.ActiveWorkbook

# ---- Products sheet -------------------------------------------------
$products = $wb.Worksheets.Item("Products")

# Header rename: product_cathegory -> product_cathegory_id
$products.Range("E1").Value = "product_cathegory_id"

# Row 6 ("Tynk akrylowy"): populate promoted_from / promoted_to (new cells),
# matching the formatting of the neighbouring E6 cell.
$products.Range("E6").Copy()
$products.Range("F6").PasteSpecial(-4122)
$products.Range("F6").Value = "2022-10-27 13:18:43.685298"

$products.Range("E6").Copy()
$products.Range("G6").PasteSpecial(-4122)
$products.Range("G6").Value = "2022-10-27 13:18:43.685298"

# Row 13 ("Syropian fundamentowy 16 cm"): add the missing F13 cathegory echo,
# matching the formatting of the neighbouring E13 cell.
$products.Range("E13").Copy()
$products.Range("F13").PasteSpecial(-4122)
$products.Range("F13").Value = "foundation_zone"

# New column widths for the newly-used F/G columns.
$products.Columns.Item(6).ColumnWidth = 22.666666666666668
$products.Columns.Item(7).ColumnWidth = 22.5

# ---- ProductsCathegories sheet ---------------------------------------
$cathegories = $wb.Worksheets.Item("ProductsCathegories")

# Remove the stray empty (but styled) B1 cell entirely.
$cathegories.Range("B1").Clear()

# ---- Users sheet -------------------------------------------------------
$users = $wb.Worksheets.Item("Users")

# Seed a default "avatars" JSON value for both users.
$users.Range("C2").Value = "{}"
$users.Range("C3").Value = "{}"
